# Auto-generated Excel COM-interop script to apply the Halicarnassus_Profits update
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (Leve Item ID 4564)
$ws.Range("H6").Value = 78.916664
$ws.Range("I6").Value = 71.14286
$ws.Range("K6").Value = 213.42858
$ws.Range("M6").Value = -101.42858

# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 650
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 650
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 650
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1302

# Row 41 (Leve Item ID 5478)
$ws.Range("H41").Value = 620.53845
$ws.Range("I41").Value = 799.5714
$ws.Range("J41").Value = 411.66666
$ws.Range("K41").Value = 799.5714
$ws.Range("L41").Value = 411.66666
$ws.Range("M41").Value = -359.5714
$ws.Range("N41").Value = -1291.66666

# Row 80 (Leve Item ID 12605)
$ws.Range("H80").Value = 542.6
$ws.Range("I80").Value = 431.125
$ws.Range("J80").Value = 670
$ws.Range("K80").Value = 1293.375
$ws.Range("L80").Value = 2010
$ws.Range("M80").Value = -295.375
$ws.Range("N80").Value = -4006

# Row 83 (Leve Item ID 12605)
$ws.Range("H83").Value = 542.6
$ws.Range("I83").Value = 431.125
$ws.Range("J83").Value = 670
$ws.Range("K83").Value = 3880.125
$ws.Range("L83").Value = 6030
$ws.Range("M83").Value = 1111.875
$ws.Range("N83").Value = -16014

# Row 94 (Leve Item ID 19905)
$ws.Range("H94").Value = 2263.75
$ws.Range("I94").Value = 2151.3635
$ws.Range("J94").Value = 3500
$ws.Range("K94").Value = 2151.3635
$ws.Range("L94").Value = 3500
$ws.Range("M94").Value = -1700.3635
$ws.Range("N94").Value = -4402

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 1754.4445
$ws.Range("I61").Value = 1298.5714
$ws.Range("J61").Value = 3350
$ws.Range("K61").Value = 1298.5714
$ws.Range("L61").Value = 3350
$ws.Range("M61").Value = -1086.5714
$ws.Range("N61").Value = -3774

# Row 74 (Leve Item ID 44000)
$ws.Range("H74").Value = 3700.8572
$ws.Range("I74").Value = 3585.4211
$ws.Range("K74").Value = 3585.4211
$ws.Range("M74").Value = -2711.4211

# Row 77 (Leve Item ID 44000)
$ws.Range("H77").Value = 3700.8572
$ws.Range("I77").Value = 3585.4211
$ws.Range("K77").Value = 17927.1055
$ws.Range("M77").Value = -13559.1055

# Row 88 (Leve Item ID 12530)
$ws.Range("H88").Value = 905.6
$ws.Range("I88").Value = 399.16666
$ws.Range("J88").Value = 1665.25
$ws.Range("K88").Value = 399.16666
$ws.Range("L88").Value = 1665.25
$ws.Range("M88").Value = 6.833340000000021
$ws.Range("N88").Value = -2477.25

# Row 91 (Leve Item ID 12530)
$ws.Range("H91").Value = 905.6
$ws.Range("I91").Value = 399.16666
$ws.Range("J91").Value = 1665.25
$ws.Range("K91").Value = 399.16666
$ws.Range("L91").Value = 1665.25
$ws.Range("M91").Value = 1004.83334
$ws.Range("N91").Value = -4473.25

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 1819.5
$ws.Range("I132").Value = 1606.55
$ws.Range("J132").Value = 3949
$ws.Range("K132").Value = 4819.65
$ws.Range("L132").Value = 11847
$ws.Range("M132").Value = -2289.65
$ws.Range("N132").Value = -16907

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 1754.4445
$ws.Range("I136").Value = 1298.5714
$ws.Range("J136").Value = 3350
$ws.Range("K136").Value = 3895.7142
$ws.Range("L136").Value = 10050
$ws.Range("M136").Value = -1345.7142
$ws.Range("N136").Value = -15150

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 7325
$ws.Range("I20").Value = 7325
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 7325
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -7078
$ws.Range("N20").ClearContents()

# Row 82 (Leve Item ID 11877)
$ws.Range("H82").Value = 11372.286
$ws.Range("I82").Value = 11372.286
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 11372.286
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -10989.286
$ws.Range("N82").ClearContents()

# Row 85 (Leve Item ID 11877)
$ws.Range("H85").Value = 11372.286
$ws.Range("I85").Value = 11372.286
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 11372.286
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -10046.286
$ws.Range("N85").ClearContents()

# Row 86 (Leve Item ID 12526)
$ws.Range("H86").Value = 1890.7222
$ws.Range("I86").Value = 1964.2307
$ws.Range("J86").Value = 1699.6
$ws.Range("K86").Value = 1964.2307
$ws.Range("L86").Value = 1699.6
$ws.Range("M86").Value = -841.2307000000001
$ws.Range("N86").Value = -3945.6

# Row 89 (Leve Item ID 12526)
$ws.Range("H89").Value = 1890.7222
$ws.Range("I89").Value = 1964.2307
$ws.Range("J89").Value = 1699.6
$ws.Range("K89").Value = 9821.1535
$ws.Range("L89").Value = 8498
$ws.Range("M89").Value = -4205.1535
$ws.Range("N89").Value = -19730

# Row 99 (Leve Item ID 19943)
$ws.Range("H99").Value = 1647.7142
$ws.Range("I99").Value = 1447
$ws.Range("K99").Value = 1447
$ws.Range("M99").Value = 51

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 5758.75
$ws.Range("I134").Value = 1561.6666
$ws.Range("K134").Value = 4684.9998
$ws.Range("M134").Value = -2149.9998

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 4504.385
$ws.Range("I31").Value = 2514.6428
$ws.Range("K31").Value = 2514.6428
$ws.Range("M31").Value = -2219.6428

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 4504.385
$ws.Range("I34").Value = 2514.6428
$ws.Range("K34").Value = 2514.6428
$ws.Range("M34").Value = -2312.6428

# Row 62 (Leve Item ID 12580)
$ws.Range("H62").Value = 2500
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 3000
$ws.Range("K62").Value = 2000
$ws.Range("L62").Value = 3000
$ws.Range("M62").Value = -1376
$ws.Range("N62").Value = -4248

# Row 65 (Leve Item ID 12580)
$ws.Range("H65").Value = 2500
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 3000
$ws.Range("K65").Value = 10000
$ws.Range("L65").Value = 15000
$ws.Range("M65").Value = -6880
$ws.Range("N65").Value = -21240

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2786.0334
$ws.Range("I132").Value = 2756.4644
$ws.Range("K132").Value = 8269.393199999999
$ws.Range("M132").Value = -5739.393199999999

$ws = $wb.Worksheets.Item("CUL")
# Row 40 (Leve Item ID 4827)
$ws.Range("H40").Value = 62.666668
$ws.Range("I40").Value = 27.416666
$ws.Range("J40").Value = 203.66667
$ws.Range("K40").Value = 109.666664
$ws.Range("L40").Value = 814.66668
$ws.Range("M40").Value = -40.666664
$ws.Range("N40").Value = -952.66668

# Row 63 (Leve Item ID 12866)
$ws.Range("H63").Value = 120
$ws.Range("I63").Value = 120
$ws.Range("K63").Value = 360
$ws.Range("M63").Value = 389

# Row 66 (Leve Item ID 12866)
$ws.Range("H66").Value = 120
$ws.Range("I66").Value = 120
$ws.Range("K66").Value = 1080
$ws.Range("M66").Value = 2664

# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 701.5
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 701.5
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

# Row 116 (Leve Item ID 27866)
$ws.Range("H116").Value = 2285.8
$ws.Range("I116").Value = 2285.8
$ws.Range("K116").Value = 6857.400000000001
$ws.Range("M116").Value = -3415.400000000001

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 2101.85
$ws.Range("I131").Value = 1485.375
$ws.Range("J131").Value = 2512.8333
$ws.Range("K131").Value = 4456.125
$ws.Range("L131").Value = 7538.499899999999
$ws.Range("M131").Value = 583.875
$ws.Range("N131").Value = -17618.4999

# Row 138 (Leve Item ID 44105)
$ws.Range("H138").Value = 4989.857
$ws.Range("J138").Value = 8266.666999999999
$ws.Range("L138").Value = 24800.001
$ws.Range("N138").Value = -35080.001

$ws = $wb.Worksheets.Item("GSM")
# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 1850
$ws.Range("I80").Value = 1850
$ws.Range("K80").Value = 1850
$ws.Range("M80").Value = -852

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 1850
$ws.Range("I83").Value = 1850
$ws.Range("K83").Value = 9250
$ws.Range("M83").Value = -4258

# Row 116 (Leve Item ID 26120)
$ws.Range("H116").Value = 48935.5
$ws.Range("J116").Value = 48935.5
$ws.Range("L116").Value = 48935.5
$ws.Range("N116").Value = -58113.5

# Row 122 (Leve Item ID 36182)
$ws.Range("H122").Value = 2560.4
$ws.Range("I122").Value = 2560.4
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 7681.200000000001
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -5231.200000000001
$ws.Range("N122").ClearContents()

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 2719.85
$ws.Range("I126").Value = 2906.125
$ws.Range("J126").Value = 1974.75
$ws.Range("K126").Value = 8718.375
$ws.Range("L126").Value = 5924.25
$ws.Range("M126").Value = -6248.375
$ws.Range("N126").Value = -10864.25

# Row 130 (Leve Item ID 34692)
$ws.Range("H130").Value = 30000
$ws.Range("J130").Value = 30000
$ws.Range("L130").Value = 30000
$ws.Range("N130").Value = -40040

$ws = $wb.Worksheets.Item("LTW")
# Row 46 (Leve Item ID 5282)
$ws.Range("H46").Value = 3003
$ws.Range("I46").Value = 3613.5
$ws.Range("J46").Value = 2758.8
$ws.Range("K46").Value = 3613.5
$ws.Range("L46").Value = 2758.8
$ws.Range("M46").Value = -3425.5
$ws.Range("N46").Value = -3134.8

# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 564.8
$ws.Range("I55").Value = 630.9524
$ws.Range("K55").Value = 630.9524
$ws.Range("M55").Value = -457.9524

# Row 68 (Leve Item ID 12563)
$ws.Range("H68").Value = 4472.3335
$ws.Range("J68").Value = 8500
$ws.Range("L68").Value = 8500
$ws.Range("N68").Value = -9998

# Row 71 (Leve Item ID 12563)
$ws.Range("H71").Value = 4472.3335
$ws.Range("J71").Value = 8500
$ws.Range("L71").Value = 42500
$ws.Range("N71").Value = -49988

$ws = $wb.Worksheets.Item("WVR")
# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 2219.8
$ws.Range("I132").Value = 2219.8
$ws.Range("K132").Value = 6659.400000000001
$ws.Range("M132").Value = -4129.400000000001

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 1816.48
$ws.Range("I136").Value = 1465.6
$ws.Range("J136").Value = 3220
$ws.Range("K136").Value = 4396.799999999999
$ws.Range("L136").Value = 9660
$ws.Range("M136").Value = -1846.799999999999
$ws.Range("N136").Value = -14760
